$d = $word.ActiveDocument

# Locate the existing "I = P/V = ... = 216mA" equation paragraph. It is the
# only OMath object in the document, so anchor off of it rather than a
# hard-coded paragraph index.
$anchorOMath = $d.Content.OMaths(1)
$anchorPara = $anchorOMath.Range.Paragraphs(1)

# Insert a new paragraph right after it. This naturally inherits the
# ListParagraph style + numPr bullet formatting from the source paragraph.
$anchorPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $anchorPara.Next()

# Seed the new paragraph with a placeholder OMath region so we get a real
# m:oMath anchor to target precisely. (Pushing InsertXML straight at the
# freshly split, otherwise-empty paragraph's own range works for the pPr,
# but leaves the lone m:oMath auto-wrapped in m:oMathPara.)
$placeholderXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><m:oMath><m:r><m:t>x</m:t></m:r></m:oMath></w:p>'
$newPara.Range.InsertXML($placeholderXml) | Out-Null

# Re-fetch the paragraph and replace just the placeholder OMath's own range
# (not the whole paragraph range) with the real equation:
#   R = V/I = 12V/216mA = 55.5Ω
# Replacing an *existing* OMath's range, instead of inserting fresh content
# into a plain paragraph range, keeps Word from auto-wrapping the result in
# an m:oMathPara -- matching the bare <m:oMath> used elsewhere in this doc.
$newPara2 = $anchorPara.Next()
$omathRange = $newPara2.Range.OMaths(1).Range
# The trailing value ends in the ohm sign (Ω, U+03A9); build that piece of
# the run text from its Unicode code point so the result is correct
# regardless of the host's console/script code page.
$omega = [string][char]0x03A9
$ohmsValue = "55.5" + $omega

$finalXml = '<m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>R</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="skw"/><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:fPr><m:num><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>V</m:t></m:r></m:num><m:den><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>I</m:t></m:r></m:den></m:f><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve">= </m:t></m:r><m:f><m:fPr><m:type m:val="skw"/><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:fPr><m:num><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>12V</m:t></m:r></m:num><m:den><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>216mA</m:t></m:r></m:den></m:f><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>' + $ohmsValue + '</m:t></m:r></m:oMath>'
$omathRange.InsertXML($finalXml) | Out-Null
